$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.220.23"
$ws.Range("E2").Value = "  +0.02%  "

$ws.Range("D3").Value = "'1.848.81"
$ws.Range("E3").Value = "  -0.35%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "'245.91"
$ws.Range("E5").Value = "  +1.90%  "

$ws.Range("D6").Value = "'0.6993"
$ws.Range("E6").Value = "  -0.08%  "

$ws.Range("D7").Value = "'1.000"
$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("D8").Value = "'0.07727"
$ws.Range("E8").Value = "  +0.16%  "

$ws.Range("D9").Value = "'0.3069"
$ws.Range("E9").Value = "  -0.70%  "

$ws.Range("D10").Value = "'23.55"
$ws.Range("E10").Value = "  -1.00%  "

$ws.Range("D11").Value = "'0.07824"
$ws.Range("E11").Value = "  +0.35%  "

$ws.Range("D12").Value = "'93.01"
$ws.Range("E12").Value = "  +0.94%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'5.136"
$ws.Range("E13").Value = "  +0.79%  "

$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "'1.848.07"
$ws.Range("E14").Value = "  -0.54%  "

$ws.Range("D15").Value = "'0.6864"
$ws.Range("E15").Value = "  -0.09%  "

$ws.Range("D16").Value = "'6.635"
$ws.Range("E16").Value = "  +2.13%  "

$ws.Range("D17").Value = "'0.000008329"
$ws.Range("E17").Value = "  -0.71%  "

$ws.Range("D18").Value = "'29.197.30"
$ws.Range("E18").Value = "  -0.07%  "

$ws.Range("D19").Value = "'242.00"
$ws.Range("E19").Value = "  -2.98%  "

$ws.Range("D20").Value = "'2.089.54"
$ws.Range("E20").Value = "  -1.10%  "

$ws.Range("D21").Value = "'12.74"
$ws.Range("E21").Value = "  -0.79%  "

$ws.Range("E22").Value = "  +0.01%  "

$ws.Range("D23").Value = "'7.522"
$ws.Range("E23").Value = "  -0.10%  "

$ws.Range("D24").Value = "'0.9999"
$ws.Range("E24").Value = "  -0.03%  "

$ws.Range("E25").Value = "  -0.72%  "

$ws.Range("D26").Value = "'159.06"
$ws.Range("E26").Value = "  -0.65%  "

$ws.Range("D27").Value = "'8.832"
$ws.Range("E27").Value = "  -0.16%  "

$ws.Range("D28").Value = "'18.29"
$ws.Range("E28").Value = "  -1.15%  "

$ws.Range("D29").Value = "'1.540"
$ws.Range("E29").Value = "  -1.28%  "

$ws.Range("D30").Value = "'4.231"
$ws.Range("E30").Value = "  -0.02%  "

$ws.Range("E31").Value = "  -0.39%  "

$ws.Range("D32").Value = "'1.204"
$ws.Range("E32").Value = "  +1.02%  "

$ws.Range("D33").Value = "'0.05124"
$ws.Range("E33").Value = "  -1.17%  "

$ws.Range("D34").Value = "'0.7957"
$ws.Range("E34").Value = "  +4.45%  "

$ws.Range("D35").Value = "'1.905"
$ws.Range("E35").Value = "  +3.36%  "

$ws.Range("E36").Value = "  -1.29%  "

$ws.Range("E37").Value = "  -0.61%  "

$ws.Range("D38").Value = "'1.327.55"
$ws.Range("E38").Value = "  +8.60%  "

$ws.Range("D39").Value = "'0.01873"
$ws.Range("E39").Value = "  +0.72%  "

$ws.Range("D40").Value = "'2.715"
$ws.Range("E40").Value = "  -0.27%  "

$ws.Range("D41").Value = "'0.9522"
$ws.Range("E41").Value = "  +6.33%  "

$ws.Range("D42").Value = "'6.065"
$ws.Range("E42").Value = "  +8.96%  "

$ws.Range("D43").Value = "'107.50"
$ws.Range("E43").Value = "  -2.00%  "

$ws.Range("E44").Value = "  +0.06%  "

$ws.Range("D45").Value = "'9.748"
$ws.Range("E45").Value = "  +2.23%  "

$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "'0.00000000123"
$ws.Range("E46").Value = "  -0.79%  "

$ws.Range("B47").Value = "RocketPoolETH"
$ws.Range("C47").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D47").Value = "'1.990.22"
$ws.Range("E47").Value = "  -1.08%  "

$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").Value = "'0.5183"
$ws.Range("E48").Value = "  +0.07%  "

$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").Value = "'64.15"
$ws.Range("E49").Value = "  -1.59%  "

$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").Value = "'1.767"
$ws.Range("E50").Value = "  +1.04%  "

$ws.Range("D51").Value = "'7.006"
$ws.Range("E51").Value = "  -0.03%  "
